# Update computed leve-profit figures (columns H-N) across all profession sheets
# to match the refreshed market-board data snapshot.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4631654.5
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 4631654.5
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 13894963.5
$ws.Range("M17").ClearContents() | Out-Null
$ws.Range("N17").Value = -13895299.5
$ws.Range("H40").Value = 4402.3335
$ws.Range("I40").Value = 4777.625
$ws.Range("J40").Value = 1400
$ws.Range("K40").Value = 4777.625
$ws.Range("L40").Value = 1400
$ws.Range("M40").Value = -4602.625
$ws.Range("N40").Value = -1750
$ws.Range("H43").Value = 13141.692
$ws.Range("I43").Value = 24319
$ws.Range("J43").Value = 6155.875
$ws.Range("K43").Value = 24319
$ws.Range("L43").Value = 6155.875
$ws.Range("M43").Value = -24250
$ws.Range("N43").Value = -6293.875
$ws.Range("H51").Value = 2621.348
$ws.Range("J51").Value = 2898.7144
$ws.Range("L51").Value = 2898.7144
$ws.Range("N51").Value = -3866.7144
$ws.Range("H70").Value = 70410.07000000001
$ws.Range("J70").Value = 95196.45
$ws.Range("L70").Value = 285589.35
$ws.Range("N70").Value = -286129.35
$ws.Range("H73").Value = 70410.07000000001
$ws.Range("J73").Value = 95196.45
$ws.Range("L73").Value = 285589.35
$ws.Range("N73").Value = -287461.35
$ws.Range("H76").Value = 3598.6
$ws.Range("I76").Value = 3598.6
$ws.Range("K76").Value = 3598.6
$ws.Range("M76").Value = -3283.6
$ws.Range("H79").Value = 3598.6
$ws.Range("I79").Value = 3598.6
$ws.Range("K79").Value = 3598.6
$ws.Range("M79").Value = -2506.6
$ws.Range("H86").Value = 4912.7144
$ws.Range("I86").Value = 3878.8
$ws.Range("K86").Value = 3878.8
$ws.Range("M86").Value = -2755.8
$ws.Range("H89").Value = 4912.7144
$ws.Range("I89").Value = 3878.8
$ws.Range("K89").Value = 19394
$ws.Range("M89").Value = -13778
$ws.Range("H95").Value = 24700
$ws.Range("J95").Value = 24700
$ws.Range("L95").Value = 24700
$ws.Range("N95").Value = -30192
$ws.Range("H98").Value = 5483
$ws.Range("I98").Value = 5248.25
$ws.Range("K98").Value = 5248.25
$ws.Range("M98").Value = -3750.25
$ws.Range("H99").Value = 437
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents() | Out-Null
$ws.Range("H112").Value = 38308.45
$ws.Range("I112").Value = 1344.8572
$ws.Range("J112").Value = 50069.59
$ws.Range("K112").Value = 4034.5716
$ws.Range("L112").Value = 150208.77
$ws.Range("M112").Value = -2926.5716
$ws.Range("N112").Value = -152424.77
$ws.Range("H121").Value = 2144.5715
$ws.Range("J121").Value = 2144.5715
$ws.Range("L121").Value = 6433.7145
$ws.Range("N121").Value = -9927.7145
$ws.Range("H122").Value = 5483
$ws.Range("I122").Value = 5248.25
$ws.Range("K122").Value = 15744.75
$ws.Range("M122").Value = -13294.75
$ws.Range("H128").Value = 110000
$ws.Range("J128").Value = 110000
$ws.Range("L128").Value = 110000
$ws.Range("N128").Value = -119960
$ws.Range("H130").Value = 105000
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 105000
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 105000
$ws.Range("M130").ClearContents() | Out-Null
$ws.Range("N130").Value = -115040
$ws.Range("H131").Value = 8460.52
$ws.Range("I131").Value = 1650.55
$ws.Range("K131").Value = 4951.65
$ws.Range("M131").Value = 88.35000000000036
$ws.Range("H132").Value = 1231.2452
$ws.Range("I132").Value = 1107.3265
$ws.Range("K132").Value = 3321.979499999999
$ws.Range("M132").Value = -791.9794999999995
$ws.Range("H135").Value = 1206.7826
$ws.Range("I135").Value = 878.4737
$ws.Range("J135").Value = 2766.25
$ws.Range("K135").Value = 7906.263300000001
$ws.Range("L135").Value = 24896.25
$ws.Range("M135").Value = -5371.263300000001
$ws.Range("N135").Value = -29966.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1133.5333
$ws.Range("I2").Value = 866.4545000000001
$ws.Range("K2").Value = 866.4545000000001
$ws.Range("M2").Value = -753.4545000000001
$ws.Range("H25").Value = 5420
$ws.Range("I25").Value = 840
$ws.Range("J25").Value = 10000
$ws.Range("K25").Value = 840
$ws.Range("L25").Value = 10000
$ws.Range("M25").Value = -438
$ws.Range("N25").Value = -10804
$ws.Range("H28").Value = 2767.75
$ws.Range("I28").Value = 2767.75
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 2767.75
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -2575.75
$ws.Range("N28").ClearContents() | Out-Null
$ws.Range("H32").Value = 6476.971
$ws.Range("I32").Value = 3687.2063
$ws.Range("J32").Value = 31584.857
$ws.Range("K32").Value = 3687.2063
$ws.Range("L32").Value = 31584.857
$ws.Range("M32").Value = -3400.2063
$ws.Range("N32").Value = -32158.857
$ws.Range("H41").Value = 8367.888999999999
$ws.Range("I41").Value = 8928.429
$ws.Range("J41").Value = 6406
$ws.Range("K41").Value = 8928.429
$ws.Range("L41").Value = 6406
$ws.Range("M41").Value = -8514.429
$ws.Range("N41").Value = -7234
$ws.Range("H45").Value = 6965.3335
$ws.Range("I45").Value = 7554.8125
$ws.Range("K45").Value = 7554.8125
$ws.Range("M45").Value = -7177.8125
$ws.Range("H61").Value = 3468.4167
$ws.Range("I61").Value = 2680.111
$ws.Range("K61").Value = 2680.111
$ws.Range("M61").Value = -2468.111
$ws.Range("H63").Value = 3572
$ws.Range("J63").Value = 3560
$ws.Range("L63").Value = 3560
$ws.Range("N63").Value = -4932
$ws.Range("H66").Value = 3572
$ws.Range("J66").Value = 3560
$ws.Range("L66").Value = 17800
$ws.Range("N66").Value = -24664
$ws.Range("H74").Value = 7579.5938
$ws.Range("I74").Value = 1151.12
$ws.Range("K74").Value = 1151.12
$ws.Range("M74").Value = -277.1199999999999
$ws.Range("H77").Value = 7579.5938
$ws.Range("I77").Value = 1151.12
$ws.Range("K77").Value = 5755.599999999999
$ws.Range("M77").Value = -1387.599999999999
$ws.Range("H99").Value = 2767.75
$ws.Range("I99").Value = 2767.75
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2767.75
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 227.25
$ws.Range("N99").ClearContents() | Out-Null
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents() | Out-Null
$ws.Range("H110").Value = 6908.227
$ws.Range("I110").Value = 7411.3125
$ws.Range("J110").Value = 5566.6665
$ws.Range("K110").Value = 7411.3125
$ws.Range("L110").Value = 5566.6665
$ws.Range("M110").Value = -5366.3125
$ws.Range("N110").Value = -9656.666499999999
$ws.Range("H116").Value = 1133.5333
$ws.Range("I116").Value = 866.4545000000001
$ws.Range("K116").Value = 866.4545000000001
$ws.Range("M116").Value = 1427.5455
$ws.Range("H122").Value = 1936.3077
$ws.Range("I122").Value = 1617.5555
$ws.Range("J122").Value = 2653.5
$ws.Range("K122").Value = 4852.666499999999
$ws.Range("L122").Value = 7960.5
$ws.Range("M122").Value = -2402.666499999999
$ws.Range("N122").Value = -12860.5
$ws.Range("H132").Value = 3725.9546
$ws.Range("I132").Value = 3570.0476
$ws.Range("K132").Value = 10710.1428
$ws.Range("M132").Value = -8180.1428
$ws.Range("H136").Value = 3468.4167
$ws.Range("I136").Value = 2680.111
$ws.Range("K136").Value = 8040.333
$ws.Range("M136").Value = -5490.333

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1133.5333
$ws.Range("I3").Value = 866.4545000000001
$ws.Range("K3").Value = 866.4545000000001
$ws.Range("M3").Value = -752.4545000000001
$ws.Range("H27").Value = 30542
$ws.Range("J27").Value = 30542
$ws.Range("L27").Value = 30542
$ws.Range("N27").Value = -30926
$ws.Range("H86").Value = 4650
$ws.Range("I86").Value = 1975
$ws.Range("J86").Value = 10000
$ws.Range("K86").Value = 1975
$ws.Range("L86").Value = 10000
$ws.Range("M86").Value = -852
$ws.Range("N86").Value = -12246
$ws.Range("H89").Value = 4650
$ws.Range("I89").Value = 1975
$ws.Range("J89").Value = 10000
$ws.Range("K89").Value = 9875
$ws.Range("L89").Value = 50000
$ws.Range("M89").Value = -4259
$ws.Range("N89").Value = -61232
$ws.Range("H107").Value = 1349.7142
$ws.Range("I107").Value = 1118.375
$ws.Range("K107").Value = 1118.375
$ws.Range("M107").Value = 801.625
$ws.Range("H134").Value = 1557.6
$ws.Range("I134").Value = 1566.7677
$ws.Range("K134").Value = 4700.3031
$ws.Range("M134").Value = -2165.3031

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 324.22223
$ws.Range("I7").Value = 391.8889
$ws.Range("J7").Value = 256.55554
$ws.Range("K7").Value = 391.8889
$ws.Range("L7").Value = 256.55554
$ws.Range("M7").Value = -278.8889
$ws.Range("N7").Value = -482.55554
$ws.Range("H31").Value = 147537
$ws.Range("I31").Value = 253622.5
$ws.Range("J31").Value = 6089.6665
$ws.Range("K31").Value = 253622.5
$ws.Range("L31").Value = 6089.6665
$ws.Range("M31").Value = -253327.5
$ws.Range("N31").Value = -6679.6665
$ws.Range("H34").Value = 147537
$ws.Range("I34").Value = 253622.5
$ws.Range("J34").Value = 6089.6665
$ws.Range("K34").Value = 253622.5
$ws.Range("L34").Value = 6089.6665
$ws.Range("M34").Value = -253420.5
$ws.Range("N34").Value = -6493.6665
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents() | Out-Null
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents() | Out-Null
$ws.Range("H69").Value = 9666.666999999999
$ws.Range("I69").Value = 9000
$ws.Range("J69").Value = 11000
$ws.Range("K69").Value = 9000
$ws.Range("L69").Value = 11000
$ws.Range("M69").Value = -8251
$ws.Range("N69").Value = -12498
$ws.Range("H72").Value = 9666.666999999999
$ws.Range("I72").Value = 9000
$ws.Range("J72").Value = 11000
$ws.Range("K72").Value = 27000
$ws.Range("L72").Value = 33000
$ws.Range("M72").Value = -23256
$ws.Range("N72").Value = -40488
$ws.Range("H97").Value = 29000
$ws.Range("J97").Value = 29000
$ws.Range("L97").Value = 29000
$ws.Range("N97").Value = -30982
$ws.Range("H99").Value = 2728.647
$ws.Range("J99").Value = 3381.75
$ws.Range("L99").Value = 3381.75
$ws.Range("N99").Value = -6377.75
$ws.Range("H105").Value = 1386.2858
$ws.Range("I105").Value = 1365.8
$ws.Range("J105").Value = 1437.5
$ws.Range("K105").Value = 1365.8
$ws.Range("L105").Value = 1437.5
$ws.Range("M105").Value = 381.2
$ws.Range("N105").Value = -4931.5
$ws.Range("H126").Value = 2728.647
$ws.Range("J126").Value = 3381.75
$ws.Range("L126").Value = 10145.25
$ws.Range("N126").Value = -15085.25
$ws.Range("H134").Value = 5955.97
$ws.Range("I134").Value = 3286.3015
$ws.Range("K134").Value = 9858.904500000001
$ws.Range("M134").Value = -7323.904500000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 784.55554
$ws.Range("I5").Value = 723
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 2169
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = -2057
$ws.Range("N5").Value = -3224
$ws.Range("H24").Value = 735.6923
$ws.Range("I24").Value = 1150
$ws.Range("J24").Value = 611.4
$ws.Range("K24").Value = 3450
$ws.Range("L24").Value = 1834.2
$ws.Range("M24").Value = -3220
$ws.Range("N24").Value = -2294.2
$ws.Range("H38").Value = 44.090908
$ws.Range("I38").Value = 36.23077
$ws.Range("J38").Value = 55.444443
$ws.Range("K38").Value = 108.69231
$ws.Range("L38").Value = 166.333329
$ws.Range("M38").Value = 238.30769
$ws.Range("N38").Value = -860.333329
$ws.Range("H82").Value = 45400
$ws.Range("I82").Value = 20000
$ws.Range("J82").Value = 51750
$ws.Range("K82").Value = 60000
$ws.Range("L82").Value = 155250
$ws.Range("M82").Value = -59594
$ws.Range("N82").Value = -156062
$ws.Range("H85").Value = 45400
$ws.Range("I85").Value = 20000
$ws.Range("J85").Value = 51750
$ws.Range("K85").Value = 60000
$ws.Range("L85").Value = 155250
$ws.Range("M85").Value = -58596
$ws.Range("N85").Value = -158058
$ws.Range("H103").Value = 528
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 528
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 1584
$ws.Range("M103").ClearContents() | Out-Null
$ws.Range("N103").Value = -3342
$ws.Range("H117").Value = 953.375
$ws.Range("I117").Value = 661.1429000000001
$ws.Range("J117").Value = 2999
$ws.Range("K117").Value = 1983.4287
$ws.Range("L117").Value = 8997
$ws.Range("M117").Value = 1458.5713
$ws.Range("N117").Value = -15881
$ws.Range("H131").Value = 22783.229
$ws.Range("I131").Value = 250475
$ws.Range("K131").Value = 751425
$ws.Range("M131").Value = -746385
$ws.Range("H135").Value = 784.55554
$ws.Range("I135").Value = 723
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 6507
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -3972
$ws.Range("N135").Value = -14070

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3119.7856
$ws.Range("I80").Value = 3064.2222
$ws.Range("K80").Value = 3064.2222
$ws.Range("M80").Value = -2066.2222
$ws.Range("H83").Value = 3119.7856
$ws.Range("I83").Value = 3064.2222
$ws.Range("K83").Value = 15321.111
$ws.Range("M83").Value = -10329.111
$ws.Range("H103").Value = 10302
$ws.Range("J103").Value = 10302
$ws.Range("L103").Value = 10302
$ws.Range("N103").Value = -12646
$ws.Range("H113").Value = 2588.1
$ws.Range("I113").Value = 2741.375
$ws.Range("K113").Value = 2741.375
$ws.Range("M113").Value = -571.375
$ws.Range("H122").Value = 4949.684
$ws.Range("I122").Value = 4382.923
$ws.Range("J122").Value = 6177.6665
$ws.Range("K122").Value = 13148.769
$ws.Range("L122").Value = 18532.9995
$ws.Range("M122").Value = -10698.769
$ws.Range("N122").Value = -23432.9995
$ws.Range("H124").Value = 29999
$ws.Range("J124").Value = 29999
$ws.Range("L124").Value = 29999
$ws.Range("N124").Value = -39819
$ws.Range("H131").Value = 133000
$ws.Range("J131").Value = 133000
$ws.Range("L131").Value = 133000
$ws.Range("N131").Value = -143080
$ws.Range("H132").Value = 2678.7112
$ws.Range("I132").Value = 2111.8484
$ws.Range("J132").Value = 4237.5835
$ws.Range("K132").Value = 6335.5452
$ws.Range("L132").Value = 12712.7505
$ws.Range("M132").Value = -3805.5452
$ws.Range("N132").Value = -17772.7505

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3516.6
$ws.Range("I22").Value = 1725
$ws.Range("J22").Value = 4711
$ws.Range("K22").Value = 1725
$ws.Range("L22").Value = 4711
$ws.Range("M22").Value = -1430
$ws.Range("N22").Value = -5301
$ws.Range("H27").Value = 3516.6
$ws.Range("I27").Value = 1725
$ws.Range("J27").Value = 4711
$ws.Range("K27").Value = 1725
$ws.Range("L27").Value = 4711
$ws.Range("M27").Value = -1618
$ws.Range("N27").Value = -4925
$ws.Range("H40").Value = 2805.889
$ws.Range("I40").Value = 2484.739
$ws.Range("K40").Value = 2484.739
$ws.Range("M40").Value = -2348.739
$ws.Range("H55").Value = 181.4
$ws.Range("I55").Value = 165.23077
$ws.Range("J55").Value = 211.42857
$ws.Range("K55").Value = 165.23077
$ws.Range("L55").Value = 211.42857
$ws.Range("M55").Value = 7.769229999999993
$ws.Range("N55").Value = -557.42857
$ws.Range("H68").Value = 372608.97
$ws.Range("I68").Value = 2345.72
$ws.Range("J68").Value = 5000899.5
$ws.Range("K68").Value = 2345.72
$ws.Range("L68").Value = 5000899.5
$ws.Range("M68").Value = -1596.72
$ws.Range("N68").Value = -5002397.5
$ws.Range("H71").Value = 372608.97
$ws.Range("I71").Value = 2345.72
$ws.Range("J71").Value = 5000899.5
$ws.Range("K71").Value = 11728.6
$ws.Range("L71").Value = 25004497.5
$ws.Range("M71").Value = -7984.599999999999
$ws.Range("N71").Value = -25011985.5
$ws.Range("H74").Value = 24999.75
$ws.Range("I74").Value = 29999.5
$ws.Range("K74").Value = 29999.5
$ws.Range("M74").Value = -29001.5
$ws.Range("H77").Value = 24999.75
$ws.Range("I77").Value = 29999.5
$ws.Range("K77").Value = 89998.5
$ws.Range("M77").Value = -85006.5
$ws.Range("H82").Value = 1796.6666
$ws.Range("J82").Value = 2334
$ws.Range("L82").Value = 2334
$ws.Range("N82").Value = -3056
$ws.Range("H85").Value = 1796.6666
$ws.Range("J85").Value = 2334
$ws.Range("L85").Value = 2334
$ws.Range("N85").Value = -4830
$ws.Range("H93").Value = 990
$ws.Range("I93").Value = 1012.65216
$ws.Range("J93").Value = 469
$ws.Range("K93").Value = 1012.65216
$ws.Range("L93").Value = 469
$ws.Range("M93").Value = 235.34784
$ws.Range("N93").Value = -2965
$ws.Range("H95").Value = 35555
$ws.Range("J95").Value = 35555
$ws.Range("L95").Value = 35555
$ws.Range("N95").Value = -41047
$ws.Range("H112").Value = 64995
$ws.Range("J112").Value = 64995
$ws.Range("L112").Value = 64995
$ws.Range("N112").Value = -67949
$ws.Range("H122").Value = 5786.727
$ws.Range("I122").Value = 5072.6665
$ws.Range("J122").Value = 9000
$ws.Range("K122").Value = 15217.9995
$ws.Range("L122").Value = 27000
$ws.Range("M122").Value = -12767.9995
$ws.Range("N122").Value = -31900
$ws.Range("H132").Value = 2865.0293
$ws.Range("I132").Value = 2281.25
$ws.Range("J132").Value = 5589.3335
$ws.Range("K132").Value = 6843.75
$ws.Range("L132").Value = 16768.0005
$ws.Range("M132").Value = -4313.75
$ws.Range("N132").Value = -21828.0005
$ws.Range("H136").Value = 3255.4783
$ws.Range("I136").Value = 2718.8
$ws.Range("K136").Value = 8156.400000000001
$ws.Range("M136").Value = -5606.400000000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 24630
$ws.Range("J28").Value = 24630
$ws.Range("L28").Value = 24630
$ws.Range("N28").Value = -25326
$ws.Range("H46").Value = 100000
$ws.Range("J46").Value = 100000
$ws.Range("L46").Value = 100000
$ws.Range("N46").Value = -100462
$ws.Range("H75").Value = 29366
$ws.Range("J75").Value = 16943.334
$ws.Range("L75").Value = 16943.334
$ws.Range("N75").Value = -18815.334
$ws.Range("H78").Value = 29366
$ws.Range("J78").Value = 16943.334
$ws.Range("L78").Value = 50830.00199999999
$ws.Range("N78").Value = -60190.00199999999
$ws.Range("H81").Value = 6279.1333
$ws.Range("I81").Value = 9493.833000000001
$ws.Range("J81").Value = 4136
$ws.Range("K81").Value = 18987.666
$ws.Range("L81").Value = 8272
$ws.Range("M81").Value = -17926.666
$ws.Range("N81").Value = -10394
$ws.Range("H84").Value = 6279.1333
$ws.Range("I84").Value = 9493.833000000001
$ws.Range("J84").Value = 4136
$ws.Range("K84").Value = 94938.33
$ws.Range("L84").Value = 41360
$ws.Range("M84").Value = -89634.33
$ws.Range("N84").Value = -51968
$ws.Range("H101").Value = 35150.25
$ws.Range("J101").Value = 35150.25
$ws.Range("L101").Value = 35150.25
$ws.Range("N101").Value = -41640.25
$ws.Range("H107").Value = 91989
$ws.Range("I107").Value = 1049.25
$ws.Range("J107").Value = 334495
$ws.Range("K107").Value = 3147.75
$ws.Range("L107").Value = 1003485
$ws.Range("M107").Value = -1227.75
$ws.Range("N107").Value = -1007325
$ws.Range("H113").Value = 2929
$ws.Range("I113").Value = 2286.182
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 6858.545999999999
$ws.Range("L113").Value = 30000
$ws.Range("M113").Value = -4688.545999999999
$ws.Range("N113").Value = -34340
$ws.Range("H122").Value = 1598.4584
$ws.Range("I122").Value = 1543.7727
$ws.Range("J122").Value = 2200
$ws.Range("K122").Value = 4631.3181
$ws.Range("L122").Value = 6600
$ws.Range("M122").Value = -2181.3181
$ws.Range("N122").Value = -11500
$ws.Range("H125").Value = 111137660
$ws.Range("I125").Value = 29000
$ws.Range("K125").Value = 29000
$ws.Range("M125").Value = -24080
$ws.Range("H132").Value = 1749.2593
$ws.Range("I132").Value = 1448.409
$ws.Range("K132").Value = 4345.227000000001
$ws.Range("M132").Value = -1815.227000000001
$ws.Range("H134").Value = 100000
$ws.Range("J134").Value = 100000
$ws.Range("L134").Value = 300000
$ws.Range("N134").Value = -305070
$ws.Range("H136").Value = 1678.2646
$ws.Range("I136").Value = 1415.5862
$ws.Range("K136").Value = 4246.7586
$ws.Range("M136").Value = -1696.7586
